# Applies the "chicken" restructuring edit described in the commit:
#  - replacing_truffle_trees_cost (row 13) gets an explicit median of "NA"
#  - truffle tree planting hours (row 14) lower/upper bumped 4->5 / 6->10
#  - chicken block (rows 58-66) is replaced by an expanded 11-row block
#    (rows 58-68) that separates parameters for two mobile-coop sizes
#    (50 vs 200 chickens) and introduces egg/feed-per-hen figures
#  - everything below the old chicken block shifts down by 2 rows
#    (handled automatically by inserting 2 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Row 13: median (column C) becomes the literal "NA"
# ---------------------------------------------------------------------
$ws.Range("C13").Value() = "NA"

# ---------------------------------------------------------------------
# 2. Row 14: lower/upper bounds updated
# ---------------------------------------------------------------------
$ws.Range("B14").Value() = 5
$ws.Range("D14").Value() = 10

# ---------------------------------------------------------------------
# 3. Insert two new rows right after the existing chicken block (after
#    row 66, before the old row 67) so the block grows from 9 to 11
#    rows (58-68) and everything below shifts from 67.. to 69..
# ---------------------------------------------------------------------
$ws.Range("A67:A68").EntireRow.Insert()

# copy the formatting (fill colour etc.) of the last "old" chicken row
# (66) onto the two freshly inserted rows so they match the block style
$ws.Range("A66:F66").Copy()
$ws.Range("A67:F68").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 4. Rewrite the whole chicken block, rows 58-68
#    columns: A variable | B lower | C median | D upper | E distribution | F label
#    Cell writes are ordered to mirror the original author's edit
#    sequence (first the new "per hen" row, then working hours, then
#    the duplicated "_1"/"_2" mobile-coop rows) so the shared-string
#    table is appended to in the same order as the canonical edit.
# ---------------------------------------------------------------------

function Set-Row($r, $a, $b, $d, $f) {
    $ws.Cells.Item($r, 1).Value() = $a        # A variable
    $ws.Cells.Item($r, 2).Value() = $b        # B lower
    $ws.Cells.Item($r, 3).Value() = "NA"      # C median
    $ws.Cells.Item($r, 4).Value() = $d        # D upper
    $ws.Cells.Item($r, 5).Value() = "posnorm" # E distribution
    $ws.Cells.Item($r, 6).Value() = $f        # F label
}

# row 58 (chicken_replacement_cost) - only B/D/F values change
Set-Row 58 "chicken_replacement_cost" 7 13 "Price per chicken"

# row 61 (feed_cost) - only B/D values change
Set-Row 61 "feed_cost" 0.28000000000000003 0.34 "Feed cost per kg"

# row 63 fully written first (introduces "egg_per_hen" / "Eggs per Hen")
Set-Row 63 "egg_per_hen" 290 310 "Eggs per Hen"

# row 64, column A only (introduces "feed_per_hen"); label filled in later
$ws.Cells.Item(64, 1).Value() = "feed_per_hen"
$ws.Cells.Item(64, 2).Value() = 30
$ws.Cells.Item(64, 3).Value() = "NA"
$ws.Cells.Item(64, 4).Value() = 40
$ws.Cells.Item(64, 5).Value() = "posnorm"

# row 62 (introduces "working_hours_chicken_1")
Set-Row 62 "working_hours_chicken_1" 250 360 "Working hours chicken"

# row 60 (introduces "maintaining_chicken_mobile_1")
Set-Row 60 "maintaining_chicken_mobile_1" 500 1500 "Maintaining Chicken mobile"

# row 59 (introduces "initial_chicken_mobile_cost_1")
Set-Row 59 "initial_chicken_mobile_cost_1" 30000 40000 "Price Chicken mobile"

# back to row 64, fill in the label (introduces "Feed per Hen")
$ws.Cells.Item(64, 6).Value() = "Feed per Hen"

# row 65 (introduces "initial_chicken_mobile_cost_2")
Set-Row 65 "initial_chicken_mobile_cost_2" 5000 8000 "Price Chicken mobile"

# row 66 (introduces "maintaining_chicken_mobile_2")
Set-Row 66 "maintaining_chicken_mobile_2" 100 500 "Maintaining Chicken mobile"

# row 67 (introduces "working_hours_chicken_2")
Set-Row 67 "working_hours_chicken_2" 60 90 "Working hours chicken"

# row 68 (eggs_price) - only B/D values change
Set-Row 68 "eggs_price" 0.25 0.3 "Price eggs"

Write-Host "done"
